$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.025342063389558
$ws.Range("D2").Value = 1.036117485227592
$ws.Range("E2").Value = 1.025706743149442
$ws.Range("F2").Value = 1.036787602862162
$ws.Range("I2").Value = 1.035345182879844
$ws.Range("J2").Value = 1.030511621666424
$ws.Range("K2").Value = 1.038912144240914
$ws.Range("L2").Value = 1.028531528041338
$ws.Range("M2").Value = 1.039580345248184
$ws.Range("N2").Value = 1.031975066370532
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.026262266472544
$ws.Range("D3").Value = 1.036660778790953
$ws.Range("E3").Value = 1.026487348017936
$ws.Range("F3").Value = 1.037940891707088
$ws.Range("I3").Value = 1.035536924204052
$ws.Range("J3").Value = 1.03107126868816
$ws.Range("K3").Value = 1.03926568017579
$ws.Range("L3").Value = 1.029119565296921
$ws.Range("M3").Value = 1.040542396958487
$ws.Range("N3").Value = 1.032535508155225
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.026858106480504
$ws.Range("D4").Value = 1.03701248814117
$ws.Range("E4").Value = 1.026993180400372
$ws.Range("F4").Value = 1.038687842016855
$ws.Range("I4").Value = 1.035659947541675
$ws.Range("J4").Value = 1.031433197402722
$ws.Range("K4").Value = 1.039493895054016
$ws.Range("L4").Value = 1.029500137862673
$ws.Range("M4").Value = 1.041165028387525
$ws.Range("N4").Value = 1.032897950850075
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.027108694258136
$ws.Range("D5").Value = 1.037160383935686
$ws.Range("E5").Value = 1.027206005600531
$ws.Range("F5").Value = 1.03900202618744
$ws.Range("I5").Value = 1.035711415576034
$ws.Range("J5").Value = 1.031585303503786
$ws.Range("K5").Value = 1.039589704686809
$ws.Range("L5").Value = 1.029660147194376
$ws.Range("M5").Value = 1.041426810974706
$ws.Range("N5").Value = 1.03305027295925
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.027150774691877
$ws.Range("D6").Value = 1.03718521842363
$ws.Range("E6").Value = 1.027241750013136
$ws.Range("F6").Value = 1.039054788836144
$ws.Range("I6").Value = 1.035720042552591
$ws.Range("J6").Value = 1.031610839907488
$ws.Range("K6").Value = 1.039605783779911
$ws.Range("L6").Value = 1.02968701441443
$ws.Range("M6").Value = 1.041470767075707
$ws.Range("N6").Value = 1.033075845627574
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.026861454467945
$ws.Range("D7").Value = 1.037014464187954
$ws.Range("E7").Value = 1.026996023501822
$ws.Range("F7").Value = 1.038692039508065
$ws.Range("I7").Value = 1.035660636246901
$ws.Range("J7").Value = 1.031435230044331
$ws.Range("K7").Value = 1.039495175786645
$ws.Range("L7").Value = 1.029502275850979
$ws.Range("M7").Value = 1.04116852622831
$ws.Range("N7").Value = 1.032899986378268
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.025652965548273
$ws.Range("D8").Value = 1.036301059297478
$ws.Range("E8").Value = 1.025970400577292
$ws.Range("F8").Value = 1.037177218137879
$ws.Range("I8").Value = 1.035410199043222
$ws.Range("J8").Value = 1.030700798013791
$ws.Range("K8").Value = 1.039031735883882
$ws.Range("L8").Value = 1.028730242308725
$ws.Range("M8").Value = 1.039905450168608
$ws.Range("N8").Value = 1.032164511370004
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.023526607468274
$ws.Range("D9").Value = 1.03504526153456
$ws.Range("E9").Value = 1.02416875765437
$ws.Range("F9").Value = 1.034513243100008
$ws.Range("I9").Value = 1.034960908211235
$ws.Range("J9").Value = 1.029405134608417
$ws.Range("K9").Value = 1.038210956940572
$ws.Range("L9").Value = 1.027370419736124
$ws.Range("M9").Value = 1.037680678509825
$ws.Range("N9").Value = 1.030867007973981
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.022111202979481
$ws.Range("D10").Value = 1.034209041852925
$ws.Range("E10").Value = 1.02297152314542
$ws.Range("F10").Value = 1.032740854699541
$ws.Range("I10").Value = 1.034656041408472
$ws.Range("J10").Value = 1.02854039519251
$ws.Range("K10").Value = 1.037661051537309
$ws.Range("L10").Value = 1.026464326018784
$ws.Range("M10").Value = 1.036198138373625
$ws.Range("N10").Value = 1.030001040528916
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.021498840779987
$ws.Range("D11").Value = 1.033847201918378
$ws.Range("E11").Value = 1.02245403797543
$ws.Range("F11").Value = 1.031974244619712
$ws.Range("I11").Value = 1.034522770742796
$ws.Range("J11").Value = 1.028165734061803
$ws.Range("K11").Value = 1.03742230464495
$ws.Range("L11").Value = 1.026072095953978
$ws.Range("M11").Value = 1.035556336762173
$ws.Range("N11").Value = 1.029625847336429
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.02127146061091
$ws.Range("D12").Value = 1.033712837453733
$ws.Range("E12").Value = 1.022261961128385
$ws.Range("F12").Value = 1.03168961854747
$ws.Range("I12").Value = 1.034473079079149
$ws.Range("J12").Value = 1.028026535394797
$ws.Range("K12").Value = 1.03733352918472
$ws.Range("L12").Value = 1.025926422274124
$ws.Range("M12").Value = 1.035317965425566
$ws.Range("N12").Value = 1.02948645099135
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.021320230869193
$ws.Range("D13").Value = 1.033741657296615
$ws.Range("E13").Value = 1.022303155901535
$ws.Range("F13").Value = 1.031750666057465
$ws.Range("I13").Value = 1.034483746654566
$ws.Range("J13").Value = 1.028056395456269
$ws.Range("K13").Value = 1.037352576069745
$ws.Range("L13").Value = 1.025957668943738
$ws.Range("M13").Value = 1.035369095872781
$ws.Range("N13").Value = 1.029516353457534
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.02148004386052
$ws.Range("D14").Value = 1.033836094504847
$ws.Range("E14").Value = 1.022438157976792
$ws.Range("F14").Value = 1.031950714748244
$ws.Range("I14").Value = 1.034518667067141
$ws.Range("J14").Value = 1.028154228522461
$ws.Range("K14").Value = 1.037414968353994
$ws.Range("L14").Value = 1.026060054149552
$ws.Range("M14").Value = 1.035536632456343
$ws.Range("N14").Value = 1.029614325457901
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.021578520331066
$ws.Range("D15").Value = 1.033894285602008
$ws.Range("E15").Value = 1.022521355822122
$ws.Range("F15").Value = 1.032073988183067
$ws.Range("I15").Value = 1.034540157655762
$ws.Range("J15").Value = 1.02821450237263
$ws.Range("K15").Value = 1.03745339784397
$ws.Range("L15").Value = 1.026123139467973
$ws.Range("M15").Value = 1.035639860241093
$ws.Range("N15").Value = 1.029674684903851
$ws.Range("B16").Value = 1.019999999999999
$ws.Range("C16").Value = 1.022151854400231
$ws.Range("D16").Value = 1.034233061357134
$ws.Range("E16").Value = 1.02300588652655
$ws.Range("F16").Value = 1.032791749903281
$ws.Range("I16").Value = 1.034664859611192
$ws.Range("J16").Value = 1.028565255573842
$ws.Range("K16").Value = 1.037676883109637
$ws.Range("L16").Value = 1.026490359519407
$ws.Range("M16").Value = 1.036240735761813
$ws.Range("N16").Value = 1.030025936214841
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.02251163069997
$ws.Range("D17").Value = 1.034445634178007
$ws.Range("E17").Value = 1.023310068387452
$ws.Range("F17").Value = 1.033242209289984
$ws.Range("I17").Value = 1.034742744509515
$ws.Range("J17").Value = 1.028785214512196
$ws.Range("K17").Value = 1.037816900455159
$ws.Range("L17").Value = 1.026720738194255
$ws.Range("M17").Value = 1.036617688818145
$ws.Range("N17").Value = 1.030246207520117
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.022721531991886
$ws.Range("D18").Value = 1.034569648195295
$ws.Range("E18").Value = 1.023487581717288
$ws.Range("F18").Value = 1.033505036237299
$ws.Range("I18").Value = 1.034788051719804
$ws.Range("J18").Value = 1.028913491251406
$ws.Range("K18").Value = 1.037898508940144
$ws.Range("L18").Value = 1.026855125169028
$ws.Range("M18").Value = 1.036837573435086
$ws.Range("N18").Value = 1.03037466642701
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.022793111344488
$ws.Range("D19").Value = 1.034611937772328
$ws.Range("E19").Value = 1.023548124268977
$ws.Range("F19").Value = 1.033594667277338
$ws.Range("I19").Value = 1.03480347964454
$ws.Range("J19").Value = 1.02895722662149
$ws.Range("K19").Value = 1.037926324901173
$ws.Range("L19").Value = 1.026900949493679
$ws.Range("M19").Value = 1.036912550825297
$ws.Range("N19").Value = 1.030418463906335
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.022473024945197
$ws.Range("D20").Value = 1.034422824643521
$ws.Range("E20").Value = 1.023277423317982
$ws.Range("F20").Value = 1.033193870819856
$ws.Range("I20").Value = 1.034734400784865
$ws.Range("J20").Value = 1.02876161722765
$ws.Range("K20").Value = 1.037801884248295
$ws.Range("L20").Value = 1.026696019599816
$ws.Range("M20").Value = 1.036577243845017
$ws.Range("N20").Value = 1.03022257672472
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.021432980747458
$ws.Range("D21").Value = 1.033808284011191
$ws.Range("E21").Value = 1.022398399351218
$ws.Range("F21").Value = 1.031891801888014
$ws.Range("I21").Value = 1.03450838908947
$ws.Range("J21").Value = 1.028125420023841
$ws.Range("K21").Value = 1.037396597974269
$ws.Range("L21").Value = 1.026029903753944
$ws.Range("M21").Value = 1.035487296490078
$ws.Range("N21").Value = 1.029585476047909
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.020779517729312
$ws.Range("D22").Value = 1.033422124358856
$ws.Range("E22").Value = 1.021846533749329
$ws.Range("F22").Value = 1.031073874505032
$ws.Range("I22").Value = 1.034365193290464
$ws.Range("J22").Value = 1.027725227929362
$ws.Range("K22").Value = 1.037141233716449
$ws.Range("L22").Value = 1.025611195011554
$ws.Range("M22").Value = 1.034802132746414
$ws.Range("N22").Value = 1.029184715634753
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.021125887565419
$ws.Range("D23").Value = 1.033626812889909
$ws.Range("E23").Value = 1.022139010819177
$ws.Range("F23").Value = 1.031507403620856
$ws.Range("I23").Value = 1.034441207556943
$ws.Range("J23").Value = 1.027937395074068
$ws.Range("K23").Value = 1.037276658421386
$ws.Range("L23").Value = 1.025833150234327
$ws.Range("M23").Value = 1.035165338703413
$ws.Range("N23").Value = 1.02939718408114
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.022490469069959
$ws.Range("D24").Value = 1.034433131216723
$ws.Range("E24").Value = 1.023292173944193
$ws.Range("F24").Value = 1.033215712641498
$ws.Range("I24").Value = 1.03473817133106
$ws.Range("J24").Value = 1.028772279891548
$ws.Range("K24").Value = 1.037808669615225
$ws.Range("L24").Value = 1.026707188834746
$ws.Range("M24").Value = 1.036595519143215
$ws.Range("N24").Value = 1.030233254530824
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.024075943416604
$ws.Range("D25").Value = 1.035369749854091
$ws.Range("E25").Value = 1.024633850289389
$ws.Range("F25").Value = 1.035201311535023
$ws.Range("I25").Value = 1.03507800410292
$ws.Range("J25").Value = 1.029740267826428
$ws.Range("K25").Value = 1.038423631905797
$ws.Range("L25").Value = 1.027721890122163
$ws.Range("M25").Value = 1.038255723515735
$ws.Range("N25").Value = 1.031202617119601
